$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$w14Ns = "xmlns:w14='http://schemas.microsoft.com/office/word/2010/wordml'"

# ------------------------------------------------------------------
# 1) "Action" header cell (row 1, col 2): add <w:spacing w:after="0"/>
#    and replace the pPr/rPr bold formatting with Times New Roman 12pt.
# ------------------------------------------------------------------
$pAction1 = $d.Paragraphs.Item(21)
$pAction1.Range.InsertXML("<w:p $wNs $w14Ns w14:paraId='2816856E' w14:textId='77777777' w:rsidR='001A41F2' w:rsidRPr='003F6B6A' w:rsidRDefault='001A41F2' w:rsidP='005C49A1'><w:pPr><w:spacing w:after='0'/><w:rPr><w:rFonts w:ascii='Times New Roman' w:hAnsi='Times New Roman'/><w:sz w:val='24'/></w:rPr></w:pPr><w:r w:rsidRPr='003F6B6A'><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Action</w:t></w:r></w:p>")

# ------------------------------------------------------------------
# 2) "{#phases} {phase_id} {phase}" paragraph -> Heading2 styled
#    "{#phases}{phase}" (phase_id placeholder removed, grammar
#    proofErr markers added in place of the old spell-check ones).
# ------------------------------------------------------------------
$pPhases = $d.Paragraphs.Item(24)
$pPhases.Range.InsertXML("<w:p $wNs $w14Ns w14:paraId='43B4E201' w14:textId='78353919' w:rsidR='001A41F2' w:rsidRDefault='00E04EC2' w:rsidP='005C49A1'><w:pPr><w:pStyle w:val='Heading2'/><w:ind w:left='406' w:hanging='396'/></w:pPr><w:r><w:t>{#</w:t></w:r><w:proofErr w:type='gramStart'/><w:r><w:t>phases}</w:t></w:r><w:r w:rsidR='00DB053D'><w:t>{</w:t></w:r><w:proofErr w:type='gramEnd'/><w:r w:rsidR='00BC06AA'><w:t>phase</w:t></w:r><w:r w:rsidR='001A41F2' w:rsidRPr='003D1511'><w:t>}</w:t></w:r></w:p>")

# ------------------------------------------------------------------
# 3) "{action_id} {action}" paragraph -> Heading3 styled "{action}"
#    (the action_id placeholder + bookmark are dropped entirely).
# ------------------------------------------------------------------
$pActionId = $d.Paragraphs.Item(28)
$pActionId.Range.InsertXML("<w:p $wNs $w14Ns w14:paraId='55294E60' w14:textId='0AB0EE11' w:rsidR='001A41F2' w:rsidRDefault='00DB053D' w:rsidP='005C49A1'><w:pPr><w:pStyle w:val='Heading3'/><w:ind w:left='1126'/></w:pPr><w:r w:rsidR='00BC06AA'><w:t>{action}</w:t></w:r></w:p>")

# ------------------------------------------------------------------
# 4) "Output: {output}{/cards}" paragraph -> collapse the "{/cards}"
#    run trio into a single "{/}" run.
# ------------------------------------------------------------------
$pOutput = $d.Paragraphs.Item(31)
$pOutput.Range.InsertXML("<w:p $wNs $w14Ns w14:paraId='58A0FC2E' w14:textId='6DC7CB85' w:rsidR='001A41F2' w:rsidRDefault='005C49A1' w:rsidP='005C49A1'><w:pPr><w:ind w:left='493'/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:lastRenderedPageBreak/><w:t>Output: {</w:t></w:r><w:r w:rsidR='00BC06AA'><w:rPr><w:b/><w:bCs/></w:rPr><w:t>output</w:t></w:r><w:r w:rsidR='001A41F2' w:rsidRPr='00A65AD0'><w:rPr><w:b/><w:bCs/></w:rPr><w:t>}</w:t></w:r><w:r w:rsidR='00DB053D'><w:rPr><w:b/><w:bCs/></w:rPr><w:br/></w:r><w:r w:rsidR='00DB053D'><w:t>{/}</w:t></w:r></w:p>")

# ------------------------------------------------------------------
# 5) "{/phases}" paragraph -> "{/}"
# ------------------------------------------------------------------
$pClosePhases = $d.Paragraphs.Item(32)
$pClosePhases.Range.InsertXML("<w:p $wNs $w14Ns w14:paraId='0CF50DDE' w14:textId='32F7685F' w:rsidR='00E04EC2' w:rsidRDefault='00E04EC2' w:rsidP='005C49A1'><w:pPr><w:ind w:left='493'/></w:pPr><w:r><w:t>{/}</w:t></w:r></w:p>")

# ------------------------------------------------------------------
# 6) Append a new (empty) trailing row to the phases/action table.
# ------------------------------------------------------------------
$t = $d.Tables(1)
$newRow = $t.Rows.Add()
$newCell1 = $t.Cell($t.Rows.Count, 1)
$newCell1.Range.InsertXML("<w:p $wNs><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr></w:p>")
$newCell2 = $t.Cell($t.Rows.Count, 2)
$newCell2.Range.InsertXML("<w:p $wNs/>")

# ------------------------------------------------------------------
# 7) Remove the stray <w:lastRenderedPageBreak/> on the paragraph
#    right after the table (just before "Appendices").
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("Appendices", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pAfterTable = $rng.Paragraphs.Item(1).Previous(1)
$pAfterTable.Range.InsertXML("<w:p $wNs $w14Ns w14:paraId='30FA7B71' w14:textId='6F9B83A2' w:rsidR='003F6B6A' w:rsidRDefault='005C49A1' w:rsidP='00650EF6'><w:r><w:br w:type='textWrapping' w:clear='all'/></w:r></w:p>")
